$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column C: "UI Design ID" ---------------------------------------
# Header (row 1) gets the same style as the existing A1/B1 header cells.
$ws.Range("C1").Value = "UI Design ID"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)   # xlPasteFormats

# Data rows (plain/default style, matching the rest of column C)
$wireframeValues = @{
    2  = "WireFrame_Reg_001"
    3  = "WireFrame_Reg_002"
    4  = "WireFrame_Reg_003"
    5  = "WireFrame_Reg_004"
    6  = "WireFrame_Reg_005"
    7  = "WireFrame_Reg_001"
    8  = "WireFrame_Reg_006"
    9  = "WireFrame_Reg_001"
    10 = "WireFrame_Reg_007"
    11 = "WireFrame_Reg_008"
    12 = "WireFrame_Reg_009"
    13 = "WireFrame_Reg_010"
    14 = "WireFrame_Reg_011"
    15 = "WireFrame_Reg_001"
    16 = "WireFrame_Reg_012"
    17 = "WireFrame_Reg_013"
    19 = "WireFrame_Reg_001"
    20 = "WireFrame_Reg_002"
}

foreach ($row in $wireframeValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $wireframeValues[$row]
}

# --- Column C formatting / sheet view ------------------------------------
$ws.Columns.Item(3).ColumnWidth = 45.6640625

$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1

$ws.Range("C17").Select()
